$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B5").Value = -0.02
$wsSummary.Range("B6").Value = 8
$wsSummary.Range("B9").Value = 50

# --- Sheet: Strategy Status ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("D4").Value = 8
$wsStatus.Range("G4").Value = 50

# --- New trade row data (Trade #8) ---
$newRow = @{
    A = 8
    B = "2026-02-17"
    C = "19:44:33"
    D = "MarketMaking"
    E = "UP"
    F = 0.03
    G = 0.03
    H = "CLOSED"
    I = 0
    J = 0
    K = 99.98999999999999
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.11
}

# --- Sheet: All Trades ---
$wsAllTrades = $wb.Worksheets.Item("All Trades")
foreach ($col in $newRow.Keys) {
    $cell = $wsAllTrades.Range("$col" + "9")
    if ($col -eq "B") {
        # "2026-02-17" looks like a date to Excel's auto-detection; force
        # it to stay plain text (matching sibling cells B2:B8) the same way
        # the other rows were authored.
        $cell.NumberFormat = "@"
        $cell.Value = $newRow[$col]
        $cell.Style = "Normal"
    } else {
        $cell.Value = $newRow[$col]
    }
}

# --- Sheet: MarketMaking ---
$wsMM = $wb.Worksheets.Item("MarketMaking")
foreach ($col in $newRow.Keys) {
    $cell = $wsMM.Range("$col" + "9")
    if ($col -eq "B") {
        $cell.NumberFormat = "@"
        $cell.Value = $newRow[$col]
        $cell.Style = "Normal"
    } else {
        $cell.Value = $newRow[$col]
    }
}
